$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.825.99"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "1.705.32"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'317.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.3944"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").Value = "'0.4052"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'1.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").Value = "'53.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "

$ws.Range("D12").Value = "'0.08908"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "'7.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").Value = "'23.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "'8.017"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Value = "1.704.29"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "'100.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").Value = "'0.07050"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'19.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").Value = "'7.090"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'14.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D24").Value = "24.812.76"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("D25").Value = "'3.216"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("D26").Value = "'2.367"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("D27").Value = "'22.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").Value = "'162.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("D29").Value = "'8.769"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.85%  "

$ws.Range("D30").Value = "'136.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").Value = "'5.181"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = "'7.669"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.98%  "

$ws.Range("D33").Value = "'0.08886"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.39%  "

$ws.Range("D34").Value = "'1.090"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.63%  "

$ws.Range("D35").Value = "'1.985"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "

$ws.Range("D36").Value = "'11.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.62%  "

$ws.Range("D37").Value = "'0.2770"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").Value = "'14.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "'0.02792"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "'0.09216"
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.7746"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.463"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'15.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.7240"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").Value = "'2.581"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("D46").Value = "'4.212"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("B48").Value = "Flow"
$ws.Range("C48").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D48").Value = "'1.333"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'141.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").Value = "'91.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.28%  "

$ws.Range("D51").Value = "'0.08001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
